$d = $word.ActiveDocument

# The document currently ends with one empty paragraph that holds only the
# "_GoBack" bookmark, immediately before the sectPr. We need to add six new
# paragraphs of text right before it (about AVANCES.docx, GENERAL.docx,
# PLAN.docx, the wording/redaction pass, COMPETENCIAS.docx, and the class
# presentation) and two new empty paragraphs right after it.
#
# Helper: create a brand-new paragraph immediately before whatever paragraph
# currently sits last in the document (i.e. right before the bookmark
# paragraph, which always stays the last paragraph throughout this loop) and
# seed it with the given text. We re-fetch Paragraphs.Last each time instead
# of reusing a cached reference, because once a paragraph object has been
# used to split the story its "identity" tracks the position it was created
# at rather than following the bookmark paragraph down the document.
function Add-ParagraphBeforeBookmark([string]$text) {
    $doc = $word.ActiveDocument
    $bookmarkPara = $doc.Paragraphs.Last
    $bookmarkPara.Range.InsertParagraphBefore()
    $newPara = $doc.Paragraphs.Item($doc.Paragraphs.Count - 1)
    $newPara.Range.InsertAfter($text)
    return $newPara
}

# Paragraph: "El resumen de lo que hicimos..." (AVANCES.docx)
Add-ParagraphBeforeBookmark("El resumen de lo que hicimos en la primera y segunda entrega está en el documento AVANCES.docx, al cual se le hicieron modificaciones de redacción.") | Out-Null

# Paragraph: "En GENERAL.docx..."
Add-ParagraphBeforeBookmark("En GENERAL.docx se encuentra lo que es el producto en sus diferentes fases. Como no se realizará la implementación nuestro documento llega hasta la parte de diseño, sin embargo, todavía no tenemos el diseño completo, nos falta agregarle links para que sea más completo y funcional, aunque ya se realizó un prototipo que está en PROTOTIPO.pdf. Todavía no podemos avanzar con el diseño porque falta que lo apruebe la maestra Laura. Los resúmenes de las dos entregas están en el documento RESUMEN.docx") | Out-Null

# Paragraph: "En la parte de PLAN.docx..." - built up from four runs, in reading order.
$planPara = Add-ParagraphBeforeBookmark("En la parte de P")
$planPara.Range.InsertAfter("LAN.docx")
$planPara.Range.InsertAfter(" se encuentra nuestro calendario, el cuál se realizó en Trello y ahí se encuentran desglosadas las actividades que se han realizado, desde las que se harán, las que se están haciendo, las que se detuvieron, las que necesitan revisión")
$planPara.Range.InsertAfter(" y las que ya están listas. También está detallado el cómo se califica la productividad de cada integrante. ")

# Paragraph: "También se hicieron modificaciones de redacción..."
Add-ParagraphBeforeBookmark("También se hicieron modificaciones de redacción y se tuvieron que redactar nuevamente algunas cosas por el cambio de producto en todos los documentos.") | Out-Null

# Paragraph: "En COMPETENCIAS.docx..."
Add-ParagraphBeforeBookmark("En COMPETENCIAS.docx se redacta el cómo se están adquiriendo las competencias. Está en un cuadro y se acopló a la segunda entrega.") | Out-Null

# Paragraph: "Agregamos la presentación..."
Add-ParagraphBeforeBookmark("Agregamos la presentación que se hizo en clase de la segunda entrega.") | Out-Null

# Two new empty paragraphs after the (untouched) bookmark paragraph.
$bookmarkPara = $d.Paragraphs.Last
$bookmarkPara.Range.InsertParagraphAfter()
$bookmarkPara.Range.InsertParagraphAfter()

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
